# Commit: "adding custom msdeploy parameter"
#
# Slide 3 hosts three screenshot diagrams stacked/side-by-side:
#   "Picture 3" (top,   rId2) - web project -> temp dir -> Replace path -> Web Deploy package
#   "Picture 4" (middle, rId3) - source DB/website content -> web deploy package -> server -> DB/web server
#   "Picture 5" (bottom, rId4) - web.config transform chain
#
# The edit removes the bottom "web.config transform" diagram (rId4) and
# enlarges the remaining two pictures so they fill the freed vertical
# space: the top diagram grows/shifts left, and the middle diagram grows
# and drops down into the bottom diagram's old slot (it keeps showing the
# DACPAC/"web deploy package" picture - only its size/position/name change).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Find the three pictures by their (pre-edit) names - shape indices would
# shift once a shape is removed, so resolve references up front.
$picture3 = $null
$picture4 = $null
$picture5 = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "Picture 3") { $picture3 = $shape }
    elseif ($shape.Name -eq "Picture 4") { $picture4 = $shape }
    elseif ($shape.Name -eq "Picture 5") { $picture5 = $shape }
}

# Drop the bottom "web.config transform" picture (rId4) entirely.
$picture5.Delete()

# Top picture ("Picture 3", rId2): grows and shifts left/down to cover
# the vacated area. (Point literals below are tuned so the Single-
# precision COM round-trip reproduces the exact target EMU coordinates:
# off=608995,2445953 ext=6565126,1551011.)
$picture3.Left = 47.95236400472443
$picture3.Top = 192.59472440944882
$picture3.Width = 516.9390564181098
$picture3.Height = 122.12685039370079

# Middle picture ("Picture 4", rId3): moves down into the freed bottom
# slot and grows to span the full width, then is renamed to "Picture 5"
# to match the resulting layout order.
# (off=608995,606138 ext=8057075,986992)
$picture4.Left = 47.95236400472443
$picture4.Top = 47.72740157480315
$picture4.Width = 634.4153748307016
$picture4.Height = 77.71590811181088
$picture4.Name = "Picture 5"
